# Automatische test-sync: 2025-08-05 20:03:50
$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append the new mail-log row (row 60) ---
$logs = $wb.Worksheets.Item("Logs")

$newRow = 60
$logs.Cells.Item($newRow, 1).Value  = "Zou jij klant Jansen nog kunnen bellen?"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #19: Zou jij klant Jansen nog kunnen bellen?"
$logs.Cells.Item($newRow, 4).Value  = "Overig"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-05 20:03:03"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Extend the conditional-formatting ranges to include the new row ---
$ranges = @("D2:D59", "G2:G59", "H2:H59", "I2:I59", "J2:J59")
foreach ($old in $ranges) {
    $col = $old.Substring(0, 1)
    $newRange = "$col" + "2:$col" + "60"
    $fc = $logs.Range($old).FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($logs.Range($newRange))
}

# --- Sheet "Dashboard": bump the "Overig" count (B3) from 12 to 13 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 13
